# "Primeiro commit do dashboard de defeitos"
#
# Core functional change: the "Days open" helper column (H) on the "Dados"
# sheet used to compute the open-duration only when BOTH the opening date
# (E) and the closing date (F) were filled in, returning 0 for defects that
# are still open (no closing date yet):
#
#   =IF(AND(E2<>"",F2<>""),F2-E2+1,0)
#
# That's misleading for still-open defects, so the formula now falls back
# to TODAY() when the defect has not been closed yet, giving the number of
# days the defect has been open so far:
#
#   =IF(E2<>"",IF(F2<>"",F2,TODAY())-E2+1,0)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados")

# Row 2 is a standalone formula; rows 3:50 are one shared-formula block
# anchored at H3 (mirrors how the sheet was already laid out).
$ws.Range("H2").Formula = '=IF(E2<>"",IF(F2<>"",F2,TODAY())-E2+1,0)'
$ws.Range("H3:H50").Formula = '=IF(E3<>"",IF(F3<>"",F3,TODAY())-E3+1,0)'

# Carry over the cursor/selection state recorded for the sheet.
$ws.Activate()
$ws.Range("K41").Select()

$wb.Save()
